$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write new data rows (69-86): lymfoid + hjerteinfarkt indicators ---

# Row 69
$ws.Range("A69").Value = "lymfoid1"
$ws.Range("B69").Value = "lymfoid"
$ws.Range("C69").Value = "Dekningsgrad: Utredning"
$ws.Range("D69").Value = "lymfoid1"
$ws.Range("E69").Value = 0.8
$ws.Range("F69").Value = 0.6
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = "Klinisk utredningsskjema for lymfoide maligniteter er viktig for å vurdere blant annet utbredelse og antatt forløp. "
$ws.Range("I69").Value = "Figuren viser dekningsgrad av primærutredning for 2018 og 2019. Informasjon om morfologi, basis for diagnosen, diagnosedato, tilhørende sykehus etc. hentes fra patologibesvarelsen der komplettheten ligger på 97,35 prosent. Stadium, prognostiske faktorer og planlagt behandling er derimot informasjon registeret kun får via utredningsmelding og registeret har derfor hatt høy fokus på å forbedre rapporteringen. Høy grad av måloppnåelse for denne kvalitetsindikatoren er 80% eller mer."

# Row 70
$ws.Range("A70").Value = "lymfoid2"
$ws.Range("B70").Value = "lymfoid"
$ws.Range("C70").Value = "Diagnostisering av non-Hodgkin lymfom"
$ws.Range("D70").Value = "lymfoid2"
$ws.Range("E70").Value = 0.95
$ws.Range("F70").Value = 0.9
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = "I følge både pakkeforløp og handlingsprogram skal alle lymfomdiagnoser stilles etter vurdering ved et universitetssykehuslaboratorium."
$ws.Range("I70").Value = "Figuren viser andelen av pasienter diagnostisert med non-Hodgkin lymfom i 2019 som har fått diagnosen stilt ved et universitetssykehus, eller ved konsultasjon fra et universitetssykehus. Høy grad av måloppnåelse for denne kvalitetsindikatoren er over 95 prosent. Noen pasienter får ikke diagnosen stilt ved hjelp av anbefalt konsultasjon med universitetssykehus med erfaring, kompetanse og nødvendig utstyr for immunhistokjemisk og molekylær lymfomdiagnostikk. Registeret vil se på de ulike årsakene til dette og oppfordrer alle laboratorier til å følge retningslinjene."

# Row 71
$ws.Range("A71").Value = "lymfoid3"
$ws.Range("B71").Value = "lymfoid"
$ws.Range("C71").Value = "Angitt stadium v/ non-Hodgkin lymfom"
$ws.Range("D71").Value = "lymfoid3"
$ws.Range("E71").Value = 0.9
$ws.Range("F71").Value = 0.8
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Stadium viser hvor utbredt sykdommen er og er en viktig prognostisk faktor."
$ws.Range("I71").Value = "Figuren viser andelen av pasienter diagnostisert med non-Hodgkin lymfom der stadium er oppgitt på utredningsmeldingen. Å angi riktig stadium er en viktig del av diagnostisering av lymfom. Dette sier noe om hvor utbredt sykdommen er. Høy grad av måloppnåelse for denne kvalitetsindikatoren er 90 prosent eller mer. Målet er nådd med god margin. Resultatet gjelder hele helseforetaket. "

# Row 72
$ws.Range("A72").Value = "lymfoid4"
$ws.Range("B72").Value = "lymfoid"
$ws.Range("C72").Value = "Flowcytometri v/ kronisk lymfatisk leukemi"
$ws.Range("D72").Value = "lymfoid4"
$ws.Range("E72").Value = 0.9
$ws.Range("F72").Value = 0.8
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = "Ved diagnostisering av kronisk lymfatisk leukemi anbefales det å bruke undersøkelsen flowcytometri."
$ws.Range("I72").Value = "Flowcytometrisk immunfenotyping er en metode for måling av fysiske og kjemiske egenskaper til enkeltceller eller partikler i væskestrøm, og immunfenotyping av lymfocytter er den anbefalte diagnostiske metoden ved mistanke om kronisk lymfatisk leukemi. Figuren viser andelen av pasienter diagnostisert med kronisk lymfatisk leukemi som har fått diagnosen påvist ved hjelp av flowcytometri fordelt på helseforetak. Høy grad av måloppnåelse for denne kvalitetsindikatoren er 90 prosent eller mer. Resultatet gjelder hele helseforetaket."

# Row 73
$ws.Range("A73").Value = "lymfoid5"
$ws.Range("B73").Value = "lymfoid"
$ws.Range("C73").Value = "Biopsi v/ myelomatose"
$ws.Range("D73").Value = "lymfoid5"
$ws.Range("E73").Value = 0.9
$ws.Range("F73").Value = 0.8
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = "Handlingsprogrammet anbefaler biopsi som undersøkelsesmetode for fastsetting av myelomatose. "
$ws.Range("I73").Value = "Figuren viser andelen pasienter diagnostisert med myelomatose i 2019 som er diagnostisert ved hjelp av biopsi (vevsprøve). Resultatet gjelder hele helseforetaket. Praksisen er noe spredt ved de ulike helseforetakene. Biopsi gir sikrere registrering via patologene, sikrere klassifisering ved vanskelig morfologi eller lymfoblastisk fenotype, sikrere tallfesting av plasmacelleandel og vil hos enkelte pasienter gi tidligere behandling. Revidert handlingsprogram for maligne blodsykdommer har endret anbefaling til å benytte seg av biopsi og utstryk, ikke enten/eller."

# Row 74
$ws.Range("A74").Value = "lymfoid6"
$ws.Range("B74").Value = "lymfoid"
$ws.Range("C74").Value = "FISH v/ myelomatose"
$ws.Range("D74").Value = "lymfoid6"
$ws.Range("E74").Value = 0.85
$ws.Range("F74").Value = 0.7
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Handlingsprogrammet anbefaler FISH som standardmetode for cytogenetisk undersøkelse ved diagnostisering av myelomatose."
$ws.Range("I74").Value = "Figuren viser andelen pasienter diagnostisert med myelomatose i 2019 som er diagnostisert ved hjelp av en FISH-analyse. FISH er forkortelse for Fluorescerens In-Situ Hybridisering og er en cytogenetisk teknikk for å detektere og lokalisere tilstedeværelse eller fravær av spesifikke DNA- sekvenser på kromosomer, altså endringer i cellenes molekylære struktur og funksjon. FISH skal være standardmetode for cytogenetisk undersøkelse ved diagnostisering av myelomatose."

# Row 75
$ws.Range("A75").Value = "kiB"
$ws.Range("B75").Value = "hjerteinfarkt"
$ws.Range("C75").Value = "Reperfusjonsbehandling ved STEMI"
$ws.Range("D75").Value = "kiB"
$ws.Range("E75").Value = 0.9
$ws.Range("F75").Value = 0.8
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = "Pasienter under 85 år som ble behandlet med blodpropp-løsende medisin eller utblokking ved mistanke om tett hjerteåre"
$ws.Range("I75").Value = "Denne kvalitetsindikatoren er definert som andel pasienter under 85 år som ble behandlet med blodproppløsende medikament eller utblokking ved alvorlig hjerteinfarkt (STEMI). De europeiske retningslinjene som Norsk cardiologisk selskap har sluttet`r`nseg til, anbefaler at pasienter med ST-elevasjonsinfarkt (STEMI) og med <12 timer fra symptomdebut`r`ntil første medisinske kontakt (FMK) blir behandlet med trombolyse og/eller koronar angiografi/PCI i`r`nsykdomsforløpet. Hos noen pasienter vil det være medisinsk korrekt å fravike anbefalingene. Høyt målnivå er ønskelig."

# Row 76
$ws.Range("A76").Value = "kiC"
$ws.Range("B76").Value = "hjerteinfarkt"
$ws.Range("C76").Value = "Reperfusjonsbehandling innen anbefalt tid ved STEMI"
$ws.Range("D76").Value = "kiC"
$ws.Range("E76").Value = 0.85
$ws.Range("F76").Value = 0.7
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = "Pasienter under 85 år som  i løpet av kort tid ble behandlet med blodpropp-løsende medisin eller utblokking  ved mistanke om tett hjerteåre"
$ws.Range("I76").Value = "Kvalitetsindikatoren er definert som andel av pasienter under 85 år som innen anbefalt tid ble behandlet med blodproppløsende medikament eller utblokking ved mistanke om tett hjerteåre og alvorlig hjerteinfarkt (STEMI). Andel pasienter under 85 år innlagt med STelevasjonsinfarkt (STEMI) og med <12 timer fra symptomdebut til første medisinske kontakt som ble`r`nbehandlet med trombolyse innen 30 minutter eller koronar angiografi/PCI innen 120 minutter etter`r`nførste medisinske kontakt. Ved STEMI skal den tette blodåren åpnes så raskt som overhode mulig.`r`nDette vil gi redusert myokardskade og redusert risiko for hjertesvikt og død. Høyt målnivå er ønskelig."

# Row 77
$ws.Range("A77").Value = "kiC1"
$ws.Range("B77").Value = "hjerteinfarkt"
$ws.Range("C77").Value = "Trombolyse innen anbefalt tid"
$ws.Range("D77").Value = "kiC1"
$ws.Range("E77").Value = 0.8
$ws.Range("F77").Value = 0.5
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = "Pasienter under 85 år som  i løpet av kort tid ble behandlet med blodpropp-løsende medisin ved mistanke om tett hjerteåre"
$ws.Range("I77").Value = " Kvalitetsindikatoren er definert som andel pasienter under 85 år med STEMI som fikk trombolyse`r`ninnen 30 minutter av de som ble behandlet med trombolyse som initial reperfusjonsmetode. De europeiske retningslinjene som Norsk cardiologisk selskap har sluttet`r`nseg til, anbefaler at pasienter som har alvorlig hjerteinfarkt med tett hjerteåre får blodproppløsende`r`nmedikament innen 30 minutter hvis pasienten ikke kan behandles med utblokking (PCI) innen 120`r`nminutter. Høyt målnivå er ønskelig."

# Row 78
$ws.Range("A78").Value = "kiC2"
$ws.Range("B78").Value = "hjerteinfarkt"
$ws.Range("C78").Value = "Primær PCI innen anbefalt tid"
$ws.Range("D78").Value = "kiC2"
$ws.Range("E78").Value = 0.85
$ws.Range("F78").Value = 0.7
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = "Pasienter under 85 år som i løpet av kort tid ble behandlet med utblokking (PCI)  ved mistanke om tett hjerteåre"
$ws.Range("I78").Value = " Andel pasienter under 85 år med STEMI og <12 timer fra symptomdebut til første`r`nmedisinske kontakt som fikk primær PCI innen 120 minutter etter første medisinske kontakt. Primær`r`nPCI er definert som angiografi eller PCI som initial metode for å oppnå reperfusjon når angiografi`r`neller PCI ble utført innen 12 timer etter første medisinske kontakt og det ikke var gitt trombolyse på`r`nforhånd. De europeiske retningslinjene som Norsk cardiologisk selskap har sluttet`r`nseg til, anbefaler at pasienter som har alvorlig hjerteinfarkt med tett hjerteåre får behandling med`r`nutblokking (PCI) innen 120 minutter. Høyt målnivå er ønskelig."

# Row 79
$ws.Range("A79").Value = "kiD"
$ws.Range("B79").Value = "hjerteinfarkt"
$ws.Range("C79").Value = "Invasivt utredet ved NSTEMI"
$ws.Range("D79").Value = "kiD"
$ws.Range("E79").Value = 0.85
$ws.Range("F79").Value = 0.7
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = "Pasienter under 85 år med mindre EKG-forandringer som ble undersøkt med røntgen av hjertets blodårer i løpet av behandlingsforløpet"
$ws.Range("I79").Value = "Kvalitetsindikatoren er definert som andel av pasienter under 85 med mindre EKG-forandringer som ble`r`nundersøkt med røntgen av hjertets blodårer i løpet av behandlingskjeden. De europeiske retningslinjene som Norsk cardiologisk selskap har sluttet`r`nseg til, anbefaler at pasienter med NSTEMI som hovedregel utredes invasivt med koronar angiografi`r`nunder sykehusoppholdet. Ved invasiv utredning kartlegges sykdomsutbredelse. Dette er av betydning`r`nfor risikostratifisering og planlegging av behandlingen. Høyt målnivå er ønskelig."

# Row 80
$ws.Range("A80").Value = "kiE"
$ws.Range("B80").Value = "hjerteinfarkt"
$ws.Range("C80").Value = "Invasivt utredet innen 72 timer ved NSTEMI"
$ws.Range("D80").Value = "kiE"
$ws.Range("E80").Value = 0.8
$ws.Range("F80").Value = 0.5
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = "Pasienter under 85 år med mindre EKG-forandringer som ble undersøkt med røntgen av hjertets blodårer innen 72 timer etter innleggelse"
$ws.Range("I80").Value = "Kvalitetsindikatoren er definert som andel pasienter under 85 år med mindre EKG-forandringer som ble`r`nundersøkt med røntgen av hjertets blodårer innen 72 timer etter innleggelse. Europeiske retningslinjer anbefaler at pasienter med NSTEMI som`r`nhovedregel gjennomgår invasiv utredning innen 24 timer etter innleggelse. Praksis varierer i Norge.`r`nFagrådet for Norsk hjerteinfarktregister har vurdert foreliggende data og funnet at det ikke foreligger`r`ngod dokumentasjon for at pasienter med NSTEMI som hovedregel bør utredes innen 24 timer.`r`nFagrådet har derfor valgt å opprettholde utredning innen 72 timer etter innleggelse som nasjonal`r`nkvalitetsindikator. Det understrekes at pasientene må risikostratifiseres, og at pasienter med høy`r`nrisiko må utredes raskt og noen umiddelbart (innen 2 timer). Høyt målnivå er ønskelig."

# Row 81
$ws.Range("A81").Value = "kiF"
$ws.Range("B81").Value = "hjerteinfarkt"
$ws.Range("C81").Value = "Utskrevet med antitrombotisk behandling"
$ws.Range("D81").Value = "kiF"
$ws.Range("E81").Value = 0.9
$ws.Range("F81").Value = 0.8
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = "Pasienter under 85 år som behandles med to medikament for å forebygge ny blodpropp etter hjerteinfarktet"
$ws.Range("I81").Value = "Kvalitetsindikatoren er definert som andel av pasienter under 85 år som behandles med to platehemmende medikament eller et platehemmende medikament og et antikoagulasjonsbehandlende medikament etter hjerteinfarktet. Antitrombotisk behandling gir prognostisk gevinst. De europeiske`r`nretningslinjene som Norsk cardiologisk selskap har sluttet seg til, anbefaler at pasientene som`r`nhovedregel behandles i 12 måneder med to medikament for å hindre blodpropp etter`r`nhjerteinfarktet. Hos noen pasienter vil det være medisinsk korrekt å fravike anbefalingene. Høyt målnivå er ønskelig."

# Row 82
$ws.Range("A82").Value = "kiG"
$ws.Range("B82").Value = "hjerteinfarkt"
$ws.Range("C82").Value = "Utskrevet med lipidsenkende medikament"
$ws.Range("D82").Value = "kiG"
$ws.Range("E82").Value = 0.9
$ws.Range("F82").Value = 0.85
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = "Pasienter under 85 år som behandles med kolesterolsenkende medisin etter hjerteinfarktet"
$ws.Range("I82").Value = " Kvalitetsindikatoren er definert som andel pasienter under 85 år som behandles med`r`nkolesterolsenkende medikament etter hjerteinfarktet. Kolesterolsenkende behandling gir prognostisk gevinst. De europeiske`r`nretningslinjene som Norsk cardiologisk selskap har sluttet seg til, anbefaler at pasientene som`r`nhovedregel behandles med lipidsenkende medikament på ubestemt tid. Høyt målnivå er ønskelig."

# Row 83
$ws.Range("A83").Value = "kiH"
$ws.Range("B83").Value = "hjerteinfarkt"
$ws.Range("C83").Value = "Ejeksjonsfraksjon (EF) målt"
$ws.Range("D83").Value = "kiH"
$ws.Range("E83").Value = 0.8
$ws.Range("F83").Value = 0.6
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = "Undersøkelse av hjertets pumpefunksjon med ultralyd"
$ws.Range("I83").Value = "Kvalitetsindikatoren er definert som andel pasienter som undersøkes med måling av hjertets`r`npumpefunksjon. De europeiske retningslinjene som Norsk cardiologisk selskap har sluttet`r`nseg til, anbefaler at ejeksjonsfraksjon (hjertets pumpefunksjon) som hovedregel blir beskrevet og`r`nmålt under sykdomsforløpet. Høyt målnivå er ønskelig."

# Row 84
$ws.Range("A84").Value = "kiJ"
$ws.Range("B84").Value = "hjerteinfarkt"
$ws.Range("C84").Value = "Utskrevet med betablokker hvis indikasjon"
$ws.Range("D84").Value = "kiJ"
$ws.Range("E84").Value = 0.85
$ws.Range("F84").Value = 0.75
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = "Pasienter under 85 år som fikk anbefalt  medisin (betablokker) for å stabilisere hjerterytmen"
$ws.Range("I84").Value = "Kvalitetsindikatoren er definert som andel pasienter under 85 år som fikk anbefalt betablokker`r`nmedikament hvis det var indikasjon for betablokker. De europeiske retningslinjene som Norsk cardiologisk selskap har sluttet`r`nseg til, anbefaler at pasientene med EF <40 %, eller som hadde hjertesvikt i tidligere sykehistorie eller`r`nsom fikk hjertesvikt i behandlingskjeden som hovedregel bør behandles med betablokker. Hos noen`r`npasienter vil det være medisinsk korrekt å fravike anbefalingene. Høyt målnivå er ønskelig."

# Row 85
$ws.Range("A85").Value = "kiK"
$ws.Range("B85").Value = "hjerteinfarkt"
$ws.Range("C85").Value = "ACE-hemmer/AII-antagonist hvis indikasjon"
$ws.Range("D85").Value = "kiK"
$ws.Range("E85").Value = 0.8
$ws.Range("F85").Value = 0.7
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = "Pasienter under 85 år med hjertesvikt eller sukkersyke som fikk anbefalt medisin (ACE/AII hemmer)"
$ws.Range("I85").Value = "Kvalitetsindikatoren er definert som andel pasienter under 85 år med hjertesvikt eller sukkersyke som`r`nfikk medikament ACE-hemmer eller AII-antagonist. De europeiske retningslinjene som Norsk cardiologisk selskap har sluttet`r`nseg til, anbefaler at pasientene med EF <40 %, eller som har hjertesvikt i tidligere sykehistorie, eller`r`nsom fikk hjertesvikt som komplikasjon i behandlingskjeden, eller som har sukkersyke, som`r`nhovedregel behandles med ACE-hemmer/AII-antagonist. Hos noen pasienter vil det være medisinsk`r`nkorrekt å fravike anbefalingene. Høyt målnivå er ønskelig."

# Row 86
$ws.Range("A86").Value = "kiA"
$ws.Range("B86").Value = "hjerteinfarkt"
$ws.Range("C86").Value = "Dekningsgrad"
$ws.Range("D86").Value = "kiA"
$ws.Range("E86").Value = 0.85
$ws.Range("F86").Value = 0.7
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = "Hvor stor andel av pasienter med hjerteinfarkt som sykehuset meldte til  Norsk hjerteinfarktregister"
$ws.Range("I86").Value = "God dekningsgrad er en forutsetning for å vurdere kvaliteten på`r`nhjerteinfarktbehandlingen ved det enkelte sykehus. Ved registrering må sykehusene rapportere`r`nhvem de behandler, hvordan de behandler og hva som oppnås ved behandlingen. På den måten blir`r`nregistrering av alle hjerteinfarkt ved et sykehus et nødvendig og viktig verktøy for kvalitetsforbedring,`r`nog en forutsetning for å kunne bedømme alle de andre kvalitetsindikatorene."

# --- Apply cell styles matching target cellXfs (fontId 18 = Calibri 11 black) ---

$style3Cells = @(
  "A75", "C75", "D75", "F75", "H75", "I75", "A76", "C76", "D76", "F76", "H76", "I76", "A77", "C77", "D77", "F77", "H77", "I77", "A78", "C78", "D78", "F78", "H78", "I78", "A79", "C79", "D79", "F79", "H79", "I79", "A80", "C80", "D80", "F80", "H80", "I80", "A81", "C81", "D81", "F81", "H81", "I81", "A82", "C82", "D82", "F82", "H82", "I82", "A83", "C83", "D83", "F83", "H83", "I83", "A84", "C84", "D84", "F84", "H84", "I84", "A85", "C85", "D85", "F85", "H85", "I85", "A86", "C86", "D86", "F86", "H86", "I86"
)
foreach ($addr in $style3Cells) {
  $ws.Range($addr).Font.Color = 0
}

$style4Cells = @(
  "E75", "E76", "E77", "E78", "E79", "E80", "E81", "E82", "E83", "E84", "E85", "E86"
)
foreach ($addr in $style4Cells) {
  $ws.Range($addr).Font.Color = 0
  $ws.Range($addr).NumberFormat = "0.00"
}

# --- Restore selection to match target sheet view ---
$ws.Range("C86").Select()
